$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Flowchart(s)" heading paragraph: append a new run ": " right
#    after the existing ")" run (before the paragraph mark).
# ------------------------------------------------------------------
$flowchartPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.TrimEnd([char]13) -eq "Flowchart(s)") {
        $flowchartPara = $cand
        break
    }
}
$fr = $flowchartPara.Range
$insertPoint = $d.Range($fr.End - 1, $fr.End - 1)
$insertPoint.InsertAfter(": ")

# ------------------------------------------------------------------
# 2) "May be included as separate pdf" -> "In a separate pdf"
# ------------------------------------------------------------------
$d.Content.Find.Execute("May be included as separate pdf", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "In a separate pdf", 2)

# ------------------------------------------------------------------
# 3) Insert the new "Git Repository" heading paragraph plus the bold
#    "Git repository evidence" / underlined URL paragraph right after
#    the "In a separate pdf" paragraph.
# ------------------------------------------------------------------
$pdfPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.TrimEnd([char]13) -eq "In a separate pdf") {
        $pdfPara = $cand
        break
    }
}
$pr = $pdfPara.Range
$afterPdf = $d.Range($pr.End - 1, $pr.End - 1)
$afterPdf.InsertParagraphAfter()
$newPara = $d.Paragraphs($pdfPara.Index + 1)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading1"/>
            </w:pPr>
            <w:r>
              <w:t>Git R</w:t>
            </w:r>
            <w:r>
              <w:t>epository</w:t>
            </w:r>
            <w:r>
              <w:t>:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>Git r</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>epository</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t xml:space="preserve"> evidence</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:u w:val="single"/>
              </w:rPr>
              <w:t>https://github.com/artrangel/Yr3_Mechatronics_Comp.git</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newPara.Range.InsertXML($xml)

# ------------------------------------------------------------------
# 4) DefaultParagraphFont style gains <w:semiHidden/>.
# ------------------------------------------------------------------
$styles = $d.Styles
$dpf = $styles | Where-Object { $_.NameLocal -eq "Default Paragraph Font" }
if ($dpf -ne $null) {
    $dpf.SemiHidden = $true
}
